# The document has a title-page section, so it carries two distinct
# header parts and two distinct footer parts (primary + first page).
# Each one contains a single inline picture (the Pearson logo in the
# footers, the BTEC logo in the headers). Word's InlineShape object
# does not expose a settable Name, so each picture is temporarily
# converted to a floating Shape (which does expose .Name), renamed,
# then converted back to an inline shape in place.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Headers (BTEC logo): rename from "image2.jpg" to "image1.jpg"
$hdrShape1 = $sec.Headers.Item(1).Range.InlineShapes.Item(1).ConvertToShape()
$hdrShape1.Name = "image1.jpg"
$hdrShape1.ConvertToInlineShape() | Out-Null

$hdrShape2 = $sec.Headers.Item(2).Range.InlineShapes.Item(1).ConvertToShape()
$hdrShape2.Name = "image1.jpg"
$hdrShape2.ConvertToInlineShape() | Out-Null

# Footers (Pearson Edexcel logo): rename from "image1.png" to "image2.png"
$ftrShape1 = $sec.Footers.Item(1).Range.InlineShapes.Item(1).ConvertToShape()
$ftrShape1.Name = "image2.png"
$ftrShape1.ConvertToInlineShape() | Out-Null

$ftrShape2 = $sec.Footers.Item(2).Range.InlineShapes.Item(1).ConvertToShape()
$ftrShape2.Name = "image2.png"
$ftrShape2.ConvertToInlineShape() | Out-Null
